$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.334194987640927
$ws.Range("C2").Value = 0.06156237666765207
$ws.Range("D2").Value = 0.1444367476148454
$ws.Range("F2").Value = 2.4946464669882
$ws.Range("G2").Value = 0.002573306288984005
$ws.Range("K2").Value = 0.8524937038761777
$ws.Range("L2").Value = 0.2456154873396983
$ws.Range("M2").Value = 0.3123384522624235
$ws.Range("N2").Value = 3.623824324771647
$ws.Range("B3").Value = 1.292403242238777
$ws.Range("C3").Value = 0.05725434632145721
$ws.Range("D3").Value = 0.1447525391802849
$ws.Range("F3").Value = 2.470508568680785
$ws.Range("G3").Value = 0.002577554362360618
$ws.Range("K3").Value = 0.809715614828292
$ws.Range("L3").Value = 0.2426147386426294
$ws.Range("M3").Value = 0.3047791175884456
$ws.Range("N3").Value = 3.625189172712837
$ws.Range("B4").Value = 1.26754056097576
$ws.Range("C4").Value = 0.05458010522372092
$ws.Range("D4").Value = 0.1449490708373498
$ws.Range("F4").Value = 2.456721591296727
$ws.Range("G4").Value = 0.002580301329484256
$ws.Range("K4").Value = 0.7839862365261467
$ws.Range("L4").Value = 0.2408924343735919
$ws.Range("M4").Value = 0.3003153948930866
$ws.Range("N4").Value = 3.626779827217547
$ws.Range("B5").Value = 1.257609308934917
$ws.Range("C5").Value = 0.05348293875961474
$ws.Range("D5").Value = 0.1450298222792457
$ws.Range("F5").Value = 2.451362999343417
$ws.Range("G5").Value = 0.002581455713847782
$ws.Range("K5").Value = 0.7736359864404108
$ws.Range("L5").Value = 0.2402208102672887
$ws.Range("M5").Value = 0.2985410813364702
$ws.Range("N5").Value = 3.627617115228119
$ws.Range("B6").Value = 1.255972339370544
$ws.Range("C6").Value = 0.05330030646504014
$ws.Range("D6").Value = 0.1450432710737326
$ws.Range("F6").Value = 2.450488888381415
$ws.Range("G6").Value = 0.002581649514095741
$ws.Range("K6").Value = 0.7719254629196257
$ws.Range("L6").Value = 0.240111113713418
$ws.Range("M6").Value = 0.2982491569879322
$ws.Range("N6").Value = 3.62776756147403
$ws.Range("B7").Value = 1.267405812947771
$ws.Range("C7").Value = 0.05456533847974754
$ws.Range("D7").Value = 0.144950157193966
$ws.Range("F7").Value = 2.456648272212064
$ws.Range("G7").Value = 0.002580316756095214
$ws.Range("K7").Value = 0.7838461044268854
$ws.Range("L7").Value = 0.2408832542034958
$ws.Range("M7").Value = 0.3002912849377779
$ws.Range("N7").Value = 3.62679035383384
$ws.Range("B8").Value = 1.319619583716047
$ws.Range("C8").Value = 0.06008295933936836
$ws.Range("D8").Value = 0.1445450862804316
$ws.Range("F8").Value = 2.486108997047111
$ws.Range("G8").Value = 0.002574742318994039
$ws.Range("K8").Value = 0.8376322997603154
$ws.Range("L8").Value = 0.2445558937129064
$ws.Range("M8").Value = 0.3096950976359878
$ws.Range("N8").Value = 3.624138607492696
$ws.Range("B9").Value = 1.428350310935571
$ws.Range("C9").Value = 0.07067599903658106
$ws.Range("D9").Value = 0.1437715868951912
$ws.Range("F9").Value = 2.552100670014823
$ws.Range("G9").Value = 0.002564905727706845
$ws.Range("K9").Value = 0.9473815945879664
$ws.Range("L9").Value = 0.2527115950446444
$ws.Range("M9").Value = 0.3295478802170422
$ws.Range("N9").Value = 3.624919874166864
$ws.Range("B10").Value = 1.512125615329865
$ws.Range("C10").Value = 0.07832624753439177
$ws.Range("D10").Value = 0.1432158866400073
$ws.Range("F10").Value = 2.605626733371764
$ws.Range("G10").Value = 0.002558339053056303
$ws.Range("K10").Value = 1.030654472356673
$ws.Range("L10").Value = 0.2592862758090035
$ws.Range("M10").Value = 0.3449988654262626
$ws.Range("N10").Value = 3.629156011125261
$ws.Range("B11").Value = 1.551088410891509
$ws.Range("C11").Value = 0.08177909456821908
$ws.Range("D11").Value = 0.1429657934685302
$ws.Range("F11").Value = 2.631079304972161
$ws.Range("G11").Value = 0.002555493542092274
$ws.Range("K11").Value = 1.069118753921089
$ws.Range("L11").Value = 0.2624041526987639
$ws.Range("M11").Value = 0.3522168759896971
$ws.Range("N11").Value = 3.631881975558429
$ws.Range("B12").Value = 1.565965570613969
$ws.Range("C12").Value = 0.08308276592853758
$ws.Range("D12").Value = 0.1428714769913553
$ws.Range("F12").Value = 2.640876628091291
$ws.Range("G12").Value = 0.002554436281153656
$ws.Range("K12").Value = 1.083768467189799
$ws.Range("L12").Value = 0.2636030868738146
$ws.Range("M12").Value = 0.3549774153447203
$ws.Range("N12").Value = 3.633029373194518
$ws.Range("B13").Value = 1.562756044509229
$ws.Range("C13").Value = 0.08280216683101571
$ws.Range("D13").Value = 0.1428917724497047
$ws.Range("F13").Value = 2.638759522889558
$ws.Range("G13").Value = 0.002554663081234535
$ws.Range("K13").Value = 1.080609643238716
$ws.Range("L13").Value = 0.2633440627969748
$ws.Range("M13").Value = 0.3543816722102875
$ws.Range("N13").Value = 3.632777135934134
$ws.Range("B14").Value = 1.552309901982596
$ws.Range("C14").Value = 0.08188642504846655
$ws.Range("D14").Value = 0.1429580262078254
$ws.Range("F14").Value = 2.631882148375013
$ws.Range("G14").Value = 0.002555406154921404
$ws.Range("K14").Value = 1.070322307006705
$ws.Range("L14").Value = 0.2625024237585052
$ws.Range("M14").Value = 0.3524434412818636
$ws.Range("N14").Value = 3.631974063749084
$ws.Range("B15").Value = 1.545927333709983
$ws.Range("C15").Value = 0.08132500808957843
$ws.Range("D15").Value = 0.1429986591426857
$ws.Range("F15").Value = 2.627690272689762
$ws.Range("G15").Value = 0.002555863945896037
$ws.Range("K15").Value = 1.06403198061733
$ws.Range("L15").Value = 0.2619892736829712
$ws.Range("M15").Value = 0.3512597662629631
$ws.Range("N15").Value = 3.631497160187024
$ws.Range("B16").Value = 1.509596480767129
$ws.Range("C16").Value = 0.0781000533270344
$ws.Range("D16").Value = 0.1432322851365422
$ws.Range("F16").Value = 2.603985577491414
$ws.Range("G16").Value = 0.002558527856400918
$ws.Range("K16").Value = 1.028152497428152
$ws.Range("L16").Value = 0.2590850708883039
$ws.Range("M16").Value = 0.3445309616908006
$ws.Range("N16").Value = 3.628993960581084
$ws.Range("B17").Value = 1.487527244458988
$ws.Range("C17").Value = 0.07611471345759924
$ws.Range("D17").Value = 0.1433762983263742
$ws.Range("F17").Value = 2.589726330878577
$ws.Range("G17").Value = 0.002560198298176607
$ws.Range("K17").Value = 1.006291125984944
$ws.Range("L17").Value = 0.2573359669547273
$ws.Range("M17").Value = 0.3404515489017541
$ws.Range("N17").Value = 3.627663128622984
$ws.Range("B18").Value = 1.474913882259045
$ws.Range("C18").Value = 0.07497022522419172
$ws.Range("D18").Value = 0.1434593847840704
$ws.Range("F18").Value = 2.581628616520604
$ws.Range("G18").Value = 0.002561172436246586
$ws.Range("K18").Value = 0.9937719035904991
$ws.Range("L18").Value = 0.2563418846869467
$ws.Range("M18").Value = 0.338122993363541
$ws.Range("N18").Value = 3.626972852647398
$ws.Range("B19").Value = 1.470656998492927
$ws.Range("C19").Value = 0.07458227686636576
$ws.Range("D19").Value = 0.1434875600231607
$ws.Range("F19").Value = 2.57890469192057
$ws.Range("G19").Value = 0.00256150455743459
$ws.Range("K19").Value = 0.9895425237098436
$ws.Range("L19").Value = 0.2560073590995415
$ws.Range("M19").Value = 0.337337643677742
$ws.Range("N19").Value = 3.62675204233021
$ws.Range("B20").Value = 1.489868242837531
$ws.Range("C20").Value = 0.07632632203663547
$ws.Range("D20").Value = 0.1433609415990968
$ws.Range("F20").Value = 2.591233503943158
$ws.Range("G20").Value = 0.002560019096449183
$ws.Range("K20").Value = 1.008612626424963
$ws.Range("L20").Value = 0.257520924677948
$ws.Range("M20").Value = 0.3408839654747808
$ws.Range("N20").Value = 3.627797015013257
$ws.Range("B21").Value = 1.555374853528008
$ws.Range("C21").Value = 0.08215550439329888
$ws.Range("D21").Value = 0.1429385553290743
$ws.Range("F21").Value = 2.633897882430688
$ws.Range("G21").Value = 0.002555187346819614
$ws.Range("K21").Value = 1.073341661601091
$ws.Range("L21").Value = 0.2627491379491289
$ws.Range("M21").Value = 0.3530120073833629
$ws.Range("N21").Value = 3.632206818797215
$ws.Range("B22").Value = 1.598903020217278
$ws.Range("C22").Value = 0.08594282565691458
$ws.Range("D22").Value = 0.1426647655778224
$ws.Range("F22").Value = 2.662708376083884
$ws.Range("G22").Value = 0.002552147634677778
$ws.Range("K22").Value = 1.116136307894948
$ws.Range("L22").Value = 0.2662725205194931
$ws.Range("M22").Value = 0.3610971417259918
$ws.Range("N22").Value = 3.635760089098653
$ws.Range("B23").Value = 1.575605674899748
$ws.Range("C23").Value = 0.08392348373072878
$ws.Range("D23").Value = 0.1428106849642408
$ws.Range("F23").Value = 2.64724675384258
$ws.Range("G23").Value = 0.002553759212687409
$ws.Range("K23").Value = 1.093251024777487
$ws.Range("L23").Value = 0.2643822864077379
$ws.Range("M23").Value = 0.3567674209329752
$ws.Range("N23").Value = 3.633802142573444
$ws.Range("B24").Value = 1.488809645314291
$ws.Range("C24").Value = 0.0762306634301666
$ws.Range("D24").Value = 0.1433678834667917
$ws.Range("F24").Value = 2.59055179916929
$ws.Range("G24").Value = 0.002560100070689852
$ws.Range("K24").Value = 1.007562922938547
$ws.Range("L24").Value = 0.2574372694684115
$ws.Range("M24").Value = 0.3406884177602123
$ws.Range("N24").Value = 3.627736251871866
$ws.Range("B25").Value = 1.398254402680124
$ws.Range("C25").Value = 0.06783399152043046
$ws.Range("D25").Value = 0.1439786202281148
$ws.Range("F25").Value = 2.533365035273874
$ws.Range("G25").Value = 0.00256745031942546
$ws.Range("K25").Value = 0.9172303645421778
$ws.Range("L25").Value = 0.2504030330087659
$ws.Range("M25").Value = 0.3240255706948645
$ws.Range("N25").Value = 3.62406656566759
